$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 187.14285
$ws.Range("I33").Value = 109.61905
$ws.Range("J33").Value = 419.7143
$ws.Range("K33").Value = 109.61905
$ws.Range("L33").Value = 419.7143
$ws.Range("M33").Value = 119.38095
$ws.Range("N33").Value = -877.7143
$ws.Range("H40").Value = 3099.875
$ws.Range("I40").Value = 1933.3334
$ws.Range("J40").Value = 3799.8
$ws.Range("K40").Value = 1933.3334
$ws.Range("L40").Value = 3799.8
$ws.Range("M40").Value = -1758.3334
$ws.Range("N40").Value = -4149.8
$ws.Range("H57").Value = 100390
$ws.Range("J57").Value = 100390
$ws.Range("L57").Value = 301170
$ws.Range("N57").Value = -302168
$ws.Range("H58").Value = 2098.2222
$ws.Range("I58").Value = 626.2857
$ws.Range("J58").Value = 7250
$ws.Range("K58").Value = 1878.8571
$ws.Range("L58").Value = 21750
$ws.Range("M58").Value = -1728.8571
$ws.Range("N58").Value = -22050
$ws.Range("H62").Value = 4975.1387
$ws.Range("I62").Value = 4742.222
$ws.Range("K62").Value = 4742.222
$ws.Range("M62").Value = -4118.222
$ws.Range("H65").Value = 4975.1387
$ws.Range("I65").Value = 4742.222
$ws.Range("K65").Value = 23711.11
$ws.Range("M65").Value = -20591.11
$ws.Range("H76").Value = 5449
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5449
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 5449
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -6079
$ws.Range("H79").Value = 5449
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5449
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 5449
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -7633
$ws.Range("H116").Value = 14694.434
$ws.Range("I116").Value = 16563
$ws.Range("K116").Value = 16563
$ws.Range("M116").Value = -13121
$ws.Range("H138").Value = 39919.37
$ws.Range("I138").Value = 2421.5625
$ws.Range("K138").Value = 7264.6875
$ws.Range("M138").Value = -2124.6875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3832.7646
$ws.Range("I2").Value = 3756.3333
$ws.Range("J2").Value = 4406
$ws.Range("K2").Value = 3756.3333
$ws.Range("L2").Value = 4406
$ws.Range("M2").Value = -3643.3333
$ws.Range("N2").Value = -4632
$ws.Range("H63").Value = 2940.9167
$ws.Range("I63").Value = 2365.889
$ws.Range("K63").Value = 2365.889
$ws.Range("M63").Value = -1679.889
$ws.Range("H66").Value = 2940.9167
$ws.Range("I66").Value = 2365.889
$ws.Range("K66").Value = 11829.445
$ws.Range("M66").Value = -8397.445
$ws.Range("H88").Value = 4738.2856
$ws.Range("I88").Value = 926.6667
$ws.Range("K88").Value = 926.6667
$ws.Range("M88").Value = -520.6667
$ws.Range("H91").Value = 4738.2856
$ws.Range("I91").Value = 926.6667
$ws.Range("K91").Value = 926.6667
$ws.Range("M91").Value = 477.3333
$ws.Range("H97").Value = 2033.2858
$ws.Range("I97").Value = 1584.375
$ws.Range("J97").Value = 2631.8333
$ws.Range("K97").Value = 1584.375
$ws.Range("L97").Value = 2631.8333
$ws.Range("M97").Value = -1088.375
$ws.Range("N97").Value = -3623.8333
$ws.Range("H102").Value = 1919.6786
$ws.Range("I102").Value = 1584.2084
$ws.Range("K102").Value = 1584.2084
$ws.Range("M102").Value = 37.79160000000002
$ws.Range("H114").Value = 49900
$ws.Range("J114").Value = 49900
$ws.Range("L114").Value = 49900
$ws.Range("N114").Value = -58578
$ws.Range("H116").Value = 3832.7646
$ws.Range("I116").Value = 3756.3333
$ws.Range("J116").Value = 4406
$ws.Range("K116").Value = 3756.3333
$ws.Range("L116").Value = 4406
$ws.Range("M116").Value = -1462.3333
$ws.Range("N116").Value = -8994
$ws.Range("H122").Value = 1825.3182
$ws.Range("I122").Value = 1607.15
$ws.Range("K122").Value = 4821.450000000001
$ws.Range("M122").Value = -2371.450000000001
$ws.Range("H132").Value = 1285.1702
$ws.Range("I132").Value = 1081.0476
$ws.Range("J132").Value = 2999.8
$ws.Range("K132").Value = 3243.142800000001
$ws.Range("L132").Value = 8999.400000000001
$ws.Range("M132").Value = -713.1428000000005
$ws.Range("N132").Value = -14059.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3832.7646
$ws.Range("I3").Value = 3756.3333
$ws.Range("J3").Value = 4406
$ws.Range("K3").Value = 3756.3333
$ws.Range("L3").Value = 4406
$ws.Range("M3").Value = -3642.3333
$ws.Range("N3").Value = -4634
$ws.Range("H7").Value = 3000
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3226
$ws.Range("H80").Value = 832.93335
$ws.Range("I80").Value = 1075.2858
$ws.Range("K80").Value = 1075.2858
$ws.Range("M80").Value = -77.28580000000011
$ws.Range("H83").Value = 832.93335
$ws.Range("I83").Value = 1075.2858
$ws.Range("K83").Value = 5376.429
$ws.Range("M83").Value = -384.4290000000001
$ws.Range("H134").Value = 1981
$ws.Range("I134").Value = 1754.3549
$ws.Range("K134").Value = 5263.0647
$ws.Range("M134").Value = -2728.0647
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2385.3333
$ws.Range("J16").Value = 3956.5
$ws.Range("L16").Value = 3956.5
$ws.Range("N16").Value = -4530.5
$ws.Range("H58").Value = 1199.0416
$ws.Range("I58").Value = 958.3158
$ws.Range("K58").Value = 958.3158
$ws.Range("M58").Value = -755.3158
$ws.Range("H60").Value = 26939.6
$ws.Range("J60").Value = 29924.5
$ws.Range("L60").Value = 29924.5
$ws.Range("N60").Value = -30946.5
$ws.Range("H94").Value = 3442.8
$ws.Range("J94").Value = 4128.6
$ws.Range("L94").Value = 4128.6
$ws.Range("N94").Value = -5030.6
$ws.Range("H105").Value = 17774.39
$ws.Range("I105").Value = 20034.95
$ws.Range("J105").Value = 2704
$ws.Range("K105").Value = 20034.95
$ws.Range("L105").Value = 2704
$ws.Range("M105").Value = -18287.95
$ws.Range("N105").Value = -6198
$ws.Range("H113").Value = 2385.3333
$ws.Range("J113").Value = 3956.5
$ws.Range("L113").Value = 3956.5
$ws.Range("N113").Value = -8296.5
$ws.Range("H136").Value = 1199.0416
$ws.Range("I136").Value = 958.3158
$ws.Range("K136").Value = 2874.9474
$ws.Range("M136").Value = -324.9474
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 198.4
$ws.Range("I12").Value = 174.57143
$ws.Range("J12").Value = 211.23077
$ws.Range("K12").Value = 523.71429
$ws.Range("L12").Value = 633.69231
$ws.Range("M12").Value = -350.71429
$ws.Range("N12").Value = -979.69231
$ws.Range("H68").Value = 4416.5557
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 4635.1763
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 13905.5289
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -15527.5289
$ws.Range("H71").Value = 4416.5557
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 4635.1763
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 41716.5867
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -49828.5867
$ws.Range("H87").Value = 16797.6
$ws.Range("I87").Value = 8747
$ws.Range("J87").Value = 49000
$ws.Range("K87").Value = 26241
$ws.Range("L87").Value = 147000
$ws.Range("M87").Value = -24993
$ws.Range("N87").Value = -149496
$ws.Range("H90").Value = 16797.6
$ws.Range("I90").Value = 8747
$ws.Range("J90").Value = 49000
$ws.Range("K90").Value = 78723
$ws.Range("L90").Value = 441000
$ws.Range("M90").Value = -72483
$ws.Range("N90").Value = -453480
$ws.Range("H124").Value = 9057.825999999999
$ws.Range("I124").Value = 2776.6667
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 8330.000100000001
$ws.Range("L124").Value = 30000
$ws.Range("M124").Value = -3420.000100000001
$ws.Range("N124").Value = -39820
$ws.Range("H131").Value = 111516.8
$ws.Range("I131").Value = 476631.22
$ws.Range("J131").Value = 1982.4667
$ws.Range("K131").Value = 1429893.66
$ws.Range("L131").Value = 5947.4001
$ws.Range("M131").Value = -1424853.66
$ws.Range("N131").Value = -16027.4001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1249.762
$ws.Range("I97").Value = 1203.2941
$ws.Range("K97").Value = 1203.2941
$ws.Range("M97").Value = -707.2941000000001
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -22940
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 84999.5
$ws.Range("I34").Value = 29999
$ws.Range("K34").Value = 29999
$ws.Range("M34").Value = -29827
$ws.Range("H46").Value = 3910.1875
$ws.Range("I46").Value = 1369.9
$ws.Range("K46").Value = 1369.9
$ws.Range("M46").Value = -1181.9
$ws.Range("H61").Value = 2298.2727
$ws.Range("I61").Value = 2348.1
$ws.Range("K61").Value = 2348.1
$ws.Range("M61").Value = -2146.1
$ws.Range("H113").Value = 2298.2727
$ws.Range("I113").Value = 2348.1
$ws.Range("K113").Value = 2348.1
$ws.Range("M113").Value = -178.0999999999999
$ws.Range("H140").Value = 127400
$ws.Range("J140").Value = 127400
$ws.Range("L140").Value = 127400
$ws.Range("N140").Value = -137760
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1454.8889
$ws.Range("I107").Value = 2050
$ws.Range("J107").Value = 2050
$ws.Range("K107").Value = 6150
$ws.Range("M107").Value = -4230

Write-Output "Applied 254 cell updates across 8 sheets"